$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.108.57"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$ws.Range("D3").Value = "1.893.20"
$ws.Range("E3").Value = "  +1.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'245.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.64%  "

# Row 6
$ws.Range("E6").Value = "  +5.84%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "'41.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.96%  "

# Row 9
$ws.Range("E9").Value = "  +4.62%  "

# Row 10
$ws.Range("D10").Value = "'52.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.54%  "

# Row 11
$ws.Range("E11").Value = "  +2.60%  "

# Row 12
$ws.Range("D12").Value = "'0.0991"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("D13").Value = "2.168.77"
$ws.Range("E13").Value = "  +1.61%  "

# Row 14
$ws.Range("D14").Value = "'12.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.76%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.903.09"
$ws.Range("E15").Value = "  +2.22%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.693"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "

# Row 17
$ws.Range("E17").Value = "  +1.05%  "

# Row 18
$ws.Range("D18").Value = "35.107.84"

# Row 19
$ws.Range("D19").Value = "'72.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.58%  "

# Row 20
$ws.Range("E20").Value = "  +1.79%  "

# Row 21
$ws.Range("D21").Value = "'239.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "

# Row 22
$ws.Range("D22").Value = "'12.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.75%  "

# Row 23
$ws.Range("E23").Value = "  +0.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +24.40%  "

# Row 26
$ws.Range("D26").Value = "'2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27
$ws.Range("E27").Value = "  +0.48%  "

# Row 29
$ws.Range("E29").Value = "  +2.85%  "

# Row 30
$ws.Range("E30").Value = "  +1.62%  "

# Row 31
$ws.Range("D31").Value = "4.166.54"
$ws.Range("E31").Value = "  +22.04%  "

# Row 32
$ws.Range("E32").Value = "  +2.06%  "

# Row 33
$ws.Range("E33").Value = "  -0.75%  "

# Row 34
$ws.Range("E34").Value = "  +14.96%  "

# Row 35
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("E36").Value = "  +0.20%  "

# Row 37
$ws.Range("E37").Value = "  -5.14%  "

# Row 38
$ws.Range("E38").Value = "  -2.43%  "

# Row 40
$ws.Range("E40").Value = "  -1.25%  "

# Row 41
$ws.Range("E41").Value = "  +2.12%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.0635"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.24%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'16.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.14%  "

# Row 44
$ws.Range("D44").Value = "'89.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.97%  "

# Row 45
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "'49.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +42.57%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.332.30"
$ws.Range("E46").Value = "  -1.21%  "

# Row 47
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("E48").Value = "  -0.30%  "

# Row 49
$ws.Range("D49").Value = "'2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.03%  "

# Row 50
$ws.Range("E50").Value = "  -2.74%  "

# Row 51
$ws.Range("D51").Value = "2.077.07"
$ws.Range("E51").Value = "  +1.29%  "
